# Updated mapping of SubjectID and RecordID
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# The "Booking Number" and "Document Control Number" mapping rows pointed at the
# wrong root element (/br-doc:BookingReport/...). Point them at the Consent
# Decision Reporting document root instead (/cdr-doc:ConsentDecisionReport/...).
$ws.Range("C15").Value = "/cdr-doc:ConsentDecisionReport/j:Booking/j:BookingSubject/j:SubjectIdentification/nc:IdentificationID"
$ws.Range("C16").Value = "cdr-doc:ConsentDecisionReport/j:Booking/j:BookingAgencyRecordIdentification/nc:IdentificationID"

# Leave the cursor on the last-edited cell, matching the author's saved selection.
$ws.Range("C16").Select()
